$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1 & 3: "Group meeting 2 (1 Hour):" and "Group meeting 3 (1 Hour):"
# headings - the " " run and "(1 Hour):" run that follow the bold/underlined
# "Group meeting N" run get combined into a single run " (1 Hour):".
# We locate each heading by its unique label text, then perform a
# case-insensitive find/replace (matching "(1 hour):" but writing back the
# correctly-cased "(1 Hour):") confined to a small range right after the
# label. This forces the engine to coalesce the two adjacent same-format
# runs that follow the label, while leaving the label run itself (and the
# unrelated "Group meeting 1" heading, which has a different run layout)
# untouched.
# ---------------------------------------------------------------------
function Fix-MeetingHourRuns($label) {
    $labelRng = $d.Range(0, $d.Content.End)
    $found = $labelRng.Find.Execute($label, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) { return }
    $afterLabel = $labelRng.End
    $tail = $d.Range($afterLabel, $afterLabel + 12)
    $tail.Find.Execute("(1 hour):", $false, $false, $false, $false, $false, $true, 0, $false, "(1 Hour):", 1)
}

Fix-MeetingHourRuns("Group meeting 2")
Fix-MeetingHourRuns("Group meeting 3")

# ---------------------------------------------------------------------
# Change 2: extend the "Spoke about the UI wire frames ..." bullet with an
# extra sentence, reproducing the exact run / proofErr split from the diff.
# We replace the whole paragraph (found via its unique existing text) with
# an equivalent paragraph carrying the same paragraph identity/formatting
# plus the additional runs, using Range.InsertXML so the proofErr markers
# around "pages" come through verbatim.
# ---------------------------------------------------------------------
$target = $d.Range(0, $d.Content.End)
$target.Find.Execute("Spoke about the UI wire frames we have designed which received positive feedback from stakeholder", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$target.Collapse(0)

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$paraXml = '<w:p ' + $wordNs + ' ' + $w14Ns + ' w14:paraId="5E2A707B" w14:textId="789D7028" w:rsidR="0065118C" w:rsidRDefault="0065118C" w:rsidP="0065118C">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
  '<w:r><w:t>Spoke about the UI wire frames we have designed which received positive feedback from stakeholder</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. He did say that it would be good if we could have a simpler design with less </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>pages</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> so we went onto redo some of the wireframe designs to better suit his needs.</w:t></w:r>' +
  '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/dummy.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $paraXml + '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

Write-Host "Done"
